$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain text (matching the source inlineStr cells), then restore
# the original (default) cell style so no stray formatting is introduced.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3815'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3414'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.44'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.194'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07478'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.447'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.068'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001094'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06649'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.622'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.410'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.482'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.533'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.006'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.069'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08703'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.654'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6921'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.432'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2206'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02342'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06319'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.779'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.236'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6495'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.854'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.143'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07136'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.89'
$ws.Range("D51").Style = "Normal"

# Remaining cells (text already, or non-numeric-looking strings) can be
# assigned directly.
$ws.Range("D2").Value = '27.245.51'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '1.785.82'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("E8").Value = '  -3.17%  '
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("E10").Value = '  -3.49%  '
$ws.Range("E11").Value = '  -3.55%  '
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("E13").Value = '  -3.19%  '
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = '1.787.29'
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("E16").Value = '  -2.23%  '
$ws.Range("E17").Value = '  -3.02%  '
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("E19").Value = '  -3.79%  '
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("E22").Value = '  -2.84%  '
$ws.Range("D23").Value = '27.264.34'
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("E24").Value = '  -6.30%  '
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E27").Value = '  -5.89%  '
$ws.Range("E28").Value = '  -3.83%  '
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("D30").Value = '1.988.76'
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("E31").Value = '  -2.07%  '
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("E33").Value = '  -5.32%  '
$ws.Range("E34").Value = '  -1.51%  '
$ws.Range("E35").Value = '  -4.72%  '
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("E38").Value = '  -4.41%  '
$ws.Range("E39").Value = '  -3.01%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E40").Value = '  -3.63%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E41").Value = '  -3.96%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E42").Value = '  -2.72%  '
$ws.Range("E43").Value = '  -4.66%  '
$ws.Range("E44").Value = '  -4.04%  '
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  -4.79%  '
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("E49").Value = '  -3.11%  '
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("E51").Value = '  -2.42%  '
